$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Инвестиции" (Investments) block added below the existing blocks,
# mirroring the layout/format of the "Сельское хозяйство" block (rows 10-12).

# Row 17: header row - category / Инвестиции (same style as row 10 header cells)
$ws.Range("B10:C10").Copy() | Out-Null
$ws.Range("B17:C17").PasteSpecial(-4122) | Out-Null

# Row 18: feature row - признаки / Инвест. в осн. кап. ... (same style as row 11)
$ws.Range("B11:C11").Copy() | Out-Null
$ws.Range("B18:C18").PasteSpecial(-4122) | Out-Null

# Row 19: trailing empty cell in column B uses the same style as B12/B5
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null

# Fill in the new cell text values (shared strings "Инвестиции" and the
# investment metric description get created automatically).
$ws.Range("B17").Value = "категория"
$ws.Range("C17").Value = "Инвестиции"
$ws.Range("B18").Value = "признаки"
$ws.Range("C18").Value = "Инвест. в осн. кап. - invest (тыс. руб) (8109001)"

# Row 19, column C is a new blank centered style distinct from B19's.
$ws.Range("C19").HorizontalAlignment = -4108

# Match the final selection recorded in the sheet view.
$ws.Range("D24").Select() | Out-Null
